# Apply updated cryptocurrency price/volume data to Sheet1
# (values mirror the authoritative commit diff for cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.299.47"
$ws.Range("E2").Value = "  -2.81%  "

$ws.Range("D3").Value = "3.000.74"
$ws.Range("E3").Value = "  -3.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.46"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.92"
$ws.Range("E6").Value = "  -7.36%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("E8").Value = "  -3.52%  "

$ws.Range("D9").Value = "2.999.72"
$ws.Range("E9").Value = "  -3.43%  "

$ws.Range("E10").Value = "  -6.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.65"
$ws.Range("E11").Value = "  -4.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -2.19%  "

$ws.Range("E13").Value = "  -5.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.60"
$ws.Range("E14").Value = "  -6.94%  "

$ws.Range("E15").Value = "  +1.81%  "

$ws.Range("D16").Value = "3.493.72"
$ws.Range("E16").Value = "  -3.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.02"
$ws.Range("E17").Value = "  -2.82%  "

$ws.Range("D18").Value = "62.290.59"
$ws.Range("E18").Value = "  -2.75%  "

$ws.Range("D19").Value = "3.000.51"
$ws.Range("E19").Value = "  -3.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "457.26"
$ws.Range("E20").Value = "  -4.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.87"
$ws.Range("E21").Value = "  -4.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.679"
$ws.Range("E22").Value = "  -5.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.31"
$ws.Range("E23").Value = "  -3.57%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.94"
$ws.Range("E24").Value = "  -1.79%  "

$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  -7.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.24"
$ws.Range("E26").Value = "  -5.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.94"
$ws.Range("E28").Value = "  -6.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").Value = "  -4.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  -3.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.09"
$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.95"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.108"
$ws.Range("E34").Value = "  -4.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.03"
$ws.Range("E35").Value = "  -3.30%  "

$ws.Range("D36").Value = "0.0₃0784"
$ws.Range("E36").Value = "  -6.97%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.74"
$ws.Range("E37").Value = "  -4.99%  "

$ws.Range("E38").Value = "  -6.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.08"
$ws.Range("E39").Value = "  -1.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.98"
$ws.Range("E40").Value = "  -2.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("E41").Value = "  -11.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "411.41"
$ws.Range("E42").Value = "  -6.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.278"
$ws.Range("E43").Value = "  -4.16%  "

$ws.Range("E44").Value = "  -1.39%  "

$ws.Range("D45").Value = "2.773.67"
$ws.Range("E45").Value = "  -2.03%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.92"
$ws.Range("E46").Value = "  -3.88%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0350"
$ws.Range("E47").Value = "  -4.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.52"
$ws.Range("E48").Value = "  -2.46%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("E50").Value = "  -2.33%  "

$ws.Range("E51").Value = "  -8.63%  "
